$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-11-27 Monday" "2023-11-28 Tuesday"

Replace-Text "94×88=8272" "68×62=4216"
Replace-Text "33×95=3135" "12×94=1128"
Replace-Text "38×47=1786" "43×23=989"
Replace-Text "16×21=336" "19×40=760"
Replace-Text "69×37=2553" "60×16=960"

Replace-Text "66×43=2838" "74×37=2738"
Replace-Text "35×87=3045" "34×79=2686"
Replace-Text "55×28=1540" "48×85=4080"
Replace-Text "52×57=2964" "75×41=3075"
Replace-Text "86×30=2580" "31×97=3007"

Replace-Text "65×21=1365" "43×75=3225"
Replace-Text "11×37=407" "90×73=6570"
Replace-Text "49×83=4067" "77×79=6083"
Replace-Text "49×91=4459" "68×55=3740"
Replace-Text "50×89=4450" "81×38=3078"

Replace-Text "29×71=2059" "89×91=8099"
Replace-Text "96×98=9408" "63×15=945"
Replace-Text "94×12=1128" "87×87=7569"
Replace-Text "40×41=1640" "89×61=5429"
Replace-Text "68×38=2584" "94×70=6580"

Replace-Text "26×59=1534" "79×23=1817"
Replace-Text "72×65=4680" "64×58=3712"
Replace-Text "30×30=900" "74×45=3330"
Replace-Text "64×12=768" "67×88=5896"
Replace-Text "35×43=1505" "14×66=924"
